# The deck currently has the "Integral" design (Red Violet colour scheme)
# applied to its single Slide Master/theme. The target state restores the
# theme's 12 scheme colours back to the default "Office" colour scheme
# (same colours PowerPoint ships with a brand-new "Office Theme" design),
# while everything else about the theme (font scheme, format scheme,
# slides, layouts, etc.) is left untouched.
#
# VBA/COM's ColorFormat.RGB takes a Long in &H00BBGGRR order (blue in the
# high byte), so build that from the familiar R,G,B bytes of each swatch.
function BGRValue($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

# Index -> scheme slot, in the standard ThemeColorScheme order:
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#  8 accent4 9 accent5 10 accent6 11 hyperlink 12 followed hyperlink
$officeColors = @(
    (BGRValue 0x00 0x00 0x00),  # dk1
    (BGRValue 0xFF 0xFF 0xFF),  # lt1
    (BGRValue 0x44 0x54 0x6A),  # dk2
    (BGRValue 0xE7 0xE6 0xE6),  # lt2
    (BGRValue 0x5B 0x9B 0xD5),  # accent1
    (BGRValue 0xED 0x7D 0x31),  # accent2
    (BGRValue 0xA5 0xA5 0xA5),  # accent3
    (BGRValue 0xFF 0xC0 0x00),  # accent4
    (BGRValue 0x44 0x72 0xC4),  # accent5
    (BGRValue 0x70 0xAD 0x47),  # accent6
    (BGRValue 0x05 0x63 0xC1),  # hlink
    (BGRValue 0x95 0x4F 0x72)   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
